$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add a "Min" column label over the Estimate column ---
$ws.Range("E3").Value = "Min"

# --- Row height bump (rows 1-4 go from 16 -> 16.5) ---
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Rows.Item(2).RowHeight = 16.5
$ws.Rows.Item(3).RowHeight = 16.5
$ws.Rows.Item(4).RowHeight = 16.5

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 21.4
$ws.Columns.Item(3).ColumnWidth = 10.17
$ws.Columns.Item(4).ColumnWidth = 10.17
$ws.Columns.Item(5).ColumnWidth = 10.17
$ws.Columns.Item(6).ColumnWidth = 12.17

# --- "Feature Block Diagram" project block (rows 5-9), entered column-first
#     like the original author (column A values, then column B values top
#     to bottom) so the shared-string table is built in the same order. ---
$ws.Range("A5").Value = "Feature Block Diagram"
$ws.Range("A12").Value = "Product Definition.txt"

$ws.Range("B5").Value = "Figure out illustrator"
$ws.Range("B6").Value = "Place Objets"
$ws.Range("B7").Value = "place relations"
$ws.Range("B8").Value = "save as pdf"
$ws.Range("B9").Value = "sync"

$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 3
$ws.Range("C8").Value = 4
$ws.Range("C9").Value = 5

$ws.Range("E5").Value = 15
$ws.Range("E6").Value = 10
$ws.Range("E7").Value = 10
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 1

$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 8
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1

# Totals row for the block above.
$ws.Range("E10").Formula = "=SUM(E5:E9)"
$ws.Range("F10").Formula = "=SUM(F5:F9)"

# D column (previously blank bordered cells) is no longer used in rows 5-6.
$ws.Range("A6").ClearFormats()
$ws.Range("A6").ClearContents()
$ws.Range("D5").ClearFormats()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearFormats()
$ws.Range("D6").ClearContents()

# --- "Product Definition.txt" project block (rows 12-13) ---
$ws.Range("B12").Value = "Fill out Text"
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 20

$ws.Range("B13").Value = "Save"
$ws.Range("C13").Value = 2
$ws.Range("E13").Value = 1

# --- Selection matches where the author ended up ---
$ws.Range("G12").Select()
